$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "nome" value in A2 (was "Brunoooooa") and the First Name
# value in E2 (was "Juliana ") to the new shared text "Julianaju".
$ws.Range("A2").Value = "Julianaju"
$ws.Range("E2").Value = "Julianaju"

# Move the active selection from A2 to E2.
[void]$ws.Range("E2").Select()
